# Generate Report for Handback
# Update the "Latest Handback DateTime" column (K) for the first data row
# (source file 3724d988-d25b-4207-bda9-92ea103e7019.md) on both the zh-cn
# and de-de localization status sheets, reflecting a newly generated
# handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-09-09 13:11:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-09 13:12:22"
